$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 30; this shifts the existing rows 30-133
# down to 31-134 and carries formatting (incl. the date style on column D).
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new weekly observation.
$ws.Range("A30").Value = 7
$ws.Range("B30").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C30").Value = "Ñuble"
$ws.Range("D30").Value = 44623
$ws.Range("E30").Value = 16
$ws.Range("F30").Value = 100112045
$ws.Range("G30").Value = "Zapallo"
$ws.Range("H30").Value = "Camote"
$ws.Range("I30").Value = "1a (cosecha)"
$ws.Range("J30").Value = 300
$ws.Range("K30").Value = 350
$ws.Range("L30").Value = 400
$ws.Range("M30").Value = 375
$ws.Range("N30").Value = "$/kilo (volumen en unidades)"
$ws.Range("O30").Value = "Región de O'Higgins"
$ws.Range("P30").Value = 375
$ws.Range("Q30").Value = 1
$ws.Range("R30").Value = "Hortaliza"
